# Add new country-specific parameters rows (23-36) to the Parameters sheet,
# using a bold Helvetica Neue 10pt font for the KEY column and a plain
# Helvetica Neue 10pt font for the VALUE column, matching the TAXDB_REGIMES-
# and-friends parameters added by this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParams = @(
    @("TAXDB_REGIMES", 6),
    @("MIN_START_YEAR", 2011),
    @("MAX_START_YEAR", 2020),
    @("MIN_START_YEAR_TRAINING", 2019),
    @("MAX_START_YEAR_TRAINING", 2019),
    @("MIN_CAPITAL_INCOME_PER_MONTH", 0),
    @("MAX_CAPITAL_INCOME_PER_MONTH", 4000),
    @("MIN_PERSONAL_PENSION_PER_MONTH", 0),
    @("MAX_PERSONAL_PENSION_PER_MONTH", 15000),
    @("MAX_CHILD_AGE_FOR_FORMAL_CARE", 14),
    @("MIN_AGE_MATERNITY", 18),
    @("MAX_AGE_MATERNITY", 44),
    @("BASE_PRICE_YEAR", 2015),
    @("PROB_NEWBORN_IS_MALE", 0.5)
)

$startRow = 23
$endRow = $startRow + $newParams.Count - 1

# Write the KEY/VALUE pairs first.
for ($i = 0; $i -lt $newParams.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newParams[$i][0]
    $ws.Cells.Item($row, 2).Value = $newParams[$i][1]
}

# Build the two distinct fonts used by these rows on a scratch cell far away
# from the used range, then fan the formatting out with copy/paste-special
# so every row shares the same two style records instead of minting a new
# one per cell.
$scratch = $ws.Cells.Item(1, 50)
$scratch.Font.Name = "Helvetica Neue"
$scratch.Font.Size = 10
# $scratch now carries the plain (non-bold) Helvetica Neue 10pt style used
# by the VALUE column.
$scratch.Copy()
$ws.Range("B$startRow`:B$endRow").PasteSpecial(-4122) | Out-Null

# Derive the bold variant used by the KEY column from the same style.
$scratch.Font.Bold = $true
$scratch.Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122) | Out-Null

$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Range("G20").Select()
